$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RiriJourney")

$ws.Range("B4").Value = 'let go with me right now babe... watch every second 😊'
$ws.Range("B6").Value = 'I''m cumming... holy fuck I''m cumming for you right now'
$ws.Range("B7").Value = 'I''m literally right on the edge babe... stay right here, I''m about to explode 🥵'
$ws.Range("B8").Value = 'oh god oh god 😊'
$ws.Range("B9").Value = 'this might be the most intense thing I''ve ever done babe... you need to see it 😊'
$ws.Range("B11").Value = 'I can feel it building so fast... you have to watch what happens next'
$ws.Range("B12").Value = 'I keep going deeper and my toes are literally curling right now 🥵'
$ws.Range("B13").Value = 'I''m fingering myself so hard right now and I can hear how wet I am babe... this is insane'
$ws.Range("B14").Value = 'FUCKK 😊'
$ws.Range("B15").Value = 'tell me you can handle this... because what I just recorded is intense 😊'
$ws.Range("B17").Value = 'be honest babe... what would you do to me right now? because I''ll act it out for you'
$ws.Range("B18").Value = 'I keep imagining you here between my legs and it''s making everything so much more intense 😊'
$ws.Range("B19").Value = 'I''m literally dripping wet right now and my hand won''t stop moving... you broke something in me'
$ws.Range("B20").Value = 'omg... okay wow that escalated 🥵'
$ws.Range("B21").Value = 'I can''t believe I''m doing this but I need you to see 😊'
$ws.Range("B23").Value = 'okay I just started touching myself and it''s 100% your fault babe... no regrets though 😊'
$ws.Range("B24").Value = 'something about this conversation is making every inch of my skin feel electric... especially between my legs'
$ws.Range("B25").Value = 'well? because I can feel my body reacting to the way you''re looking at me right now 🥵'
